$d = $word.ActiveDocument

# --- Phase 1: build all 5 paragraphs from the single pristine paragraph, ---
# --- before any direct formatting is applied, so later-created paragraphs ---
# --- don't inherit indent/highlight from already-formatted neighbours.   ---

$p1 = $d.Paragraphs(1)
$p1.Range.Text = "Files\\2018 Case Study\\CS3_Primary Sources_Policy_Strategies\\2017 National Security Strategy - § 2 references coded [ 0.06% Coverage]"
$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs(2)
$p2.Range.Text = "Reference 1 - 0.03% Coverage"
$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs(3)
$p3.Range.Text = "Malicious activity must be defeated within a network and not be passed on to its destination whenever possible."
$p3.Range.InsertParagraphAfter()

$p4 = $d.Paragraphs(4)
$p4.Range.Text = "Reference 2 - 0.02% Coverage"
$p4.Range.InsertParagraphAfter()

$p5 = $d.Paragraphs(5)
$p5.Range.Text = "This will require a resilient forward posture and agile global mobility forces."

Write-Host "Paragraph count after build:" $d.Paragraphs.Count

# --- Phase 2: apply "quote header" formatting (highlight + tight spacing + ---
# --- small indent) to paragraphs 1, 2 and 4 only.                         ---

$p1 = $d.Paragraphs(1)
$p1.Range.HighlightColorIndex = 16
$p1.SpaceBefore = 5.65
$p1.SpaceAfter = 5.65
$p1.LeftIndent = 5.65
$p1.RightIndent = 5.65
$p1.FirstLineIndent = -0.001

$p2 = $d.Paragraphs(2)
$p2.Range.HighlightColorIndex = 16
$p2.SpaceBefore = 5.65
$p2.SpaceAfter = 5.65
$p2.LeftIndent = 5.65
$p2.RightIndent = 5.65
$p2.FirstLineIndent = -0.001

$p4 = $d.Paragraphs(4)
$p4.Range.HighlightColorIndex = 16
$p4.SpaceBefore = 5.65
$p4.SpaceAfter = 5.65
$p4.LeftIndent = 5.65
$p4.RightIndent = 5.65
$p4.FirstLineIndent = -0.001

# --- Phase 3: the quoted-text paragraphs (3 and 5) keep the TextBody ---
# --- defaults except spacing, which collapses to 0/0.                ---

$p3 = $d.Paragraphs(3)
$p3.SpaceBefore = 0
$p3.SpaceAfter = 0

$p5 = $d.Paragraphs(5)
$p5.SpaceBefore = 0
$p5.SpaceAfter = 0

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
